# Update marksheet totals ("Corr/total marks" fix) on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" -> Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right column (B12): 81 -> 135
$ws.Range("B12").Value = 135

# Row 12 "Total" -> Max column (E12): "80/84" -> "135/140"
$ws.Range("E12").Value = "135/140"
